# Preliminary full sensitivity analysis: add the "mads_efast_tightened" sheet
# (tightened bounds derived from the mads_efast eFAST results), and refresh
# the stale selection/active-tab bookkeeping left over from editing.

$wb = $excel.ActiveWorkbook

# --- update stale cell selections on the existing sheets (left over from
#     where the author's cursor last was before saving) -------------------
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("D14").Select()

$wsMadsTightened = $wb.Worksheets.Item("mads_tightened")
$wsMadsTightened.Range("D25").Select()

$wsMadsEfast = $wb.Worksheets.Item("mads_efast")
$wsMadsEfast.Range("C2").Select()

# --- add the new sheet right after "mads_efast" ---------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMadsEfast)
$newSheet.Name = "mads_efast_tightened"

$newSheet.Columns.Item(1).ColumnWidth = 21.6564625850340
$newSheet.Columns.Item(2).ColumnWidth = 10.6870748299320

# Header row mirrors mads_efast's header row.
$newSheet.Range("A1").Formula = "=mads_efast!A1"
$newSheet.Range("B1").Formula = "=mads_efast!B1"
$newSheet.Range("C1").Formula = "=mads_efast!C1"
$newSheet.Range("D1").Formula = "=mads_efast!D1"

$paramRows = 2..12

foreach ($r in $paramRows) {
    $newSheet.Range("A$r").Formula = "=mads_efast!A$r"
    $newSheet.Range("B$r").Formula = "=mads_efast!B$r"

    if ($r -eq 5) {
        # "fraction" keeps the min/max pulled straight from mads_efast
        # rather than re-deriving them from the (rounded) init value.
        $newSheet.Range("C$r").Formula = "=mads_efast!C$r"
        $newSheet.Range("D$r").Formula = "=mads_efast!D$r"
    } else {
        $newSheet.Range("C$r").Formula = "=B$r-0.5"
        $newSheet.Range("D$r").Formula = "=B$r+0.5"
    }
}

# Highlight the upper bound of "fraction" in red, since it is pinned at the
# physical ceiling (0) instead of following the usual +/-0.5 tightening.
# (BGR-packed long, same encoding VBA's RGB(255,0,0) produces.)
$newSheet.Range("D5").Font.Color = 255

# Select the new sheet's own last-active cell and make it the active tab,
# matching what becomes the workbook's active sheet.
$newSheet.Range("B26").Select()

# Best-effort cosmetic restore of the tab-bar/scrollbar split the author had
# when they last saved (no functional effect on the data).
$excel.ActiveWindow.TabRatio = 419
